$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'hockey knee pads'
$ws.Cells.Item(2, 1).Value = 'knee black leggings'
$ws.Cells.Item(3, 1).Value = 'knee compression leggings'
$ws.Cells.Item(4, 1).Value = 'knee compression pad'
$ws.Cells.Item(5, 1).Value = 'knee compression running'
$ws.Cells.Item(6, 1).Value = 'knee gel pad'
$ws.Cells.Item(7, 1).Value = 'knee pad for men'
$ws.Cells.Item(8, 1).Value = 'knee pad for volleyball'
$ws.Cells.Item(9, 1).Value = 'knee pad for workout'
$ws.Cells.Item(10, 1).Value = 'knee pad leggings'
$ws.Cells.Item(11, 1).Value = 'knee pads athletic'
$ws.Cells.Item(12, 1).Value = 'knee pads basketball youth'
$ws.Cells.Item(13, 1).Value = 'knee pads black'
$ws.Cells.Item(14, 1).Value = 'knee pads boys'
$ws.Cells.Item(15, 1).Value = 'knee pads extra large'
$ws.Cells.Item(16, 1).Value = 'knee pads for football'
$ws.Cells.Item(17, 1).Value = 'knee pads for gym men'
$ws.Cells.Item(18, 1).Value = 'knee pads for men'
$ws.Cells.Item(19, 1).Value = 'knee pads for running men'
$ws.Cells.Item(20, 1).Value = 'knee pads for soccer'
$ws.Cells.Item(21, 1).Value = 'knee pads honeycomb'
$ws.Cells.Item(22, 1).Value = 'knee pads mens'
$ws.Cells.Item(23, 1).Value = 'knee pads skating youth'
$ws.Cells.Item(24, 1).Value = 'knee pads snowboarding'
$ws.Cells.Item(25, 1).Value = 'knee pads soccer'
$ws.Cells.Item(26, 1).Value = 'knee pads sport'
$ws.Cells.Item(27, 1).Value = 'knee pads squat'
$ws.Cells.Item(28, 1).Value = 'knee pads thick'
$ws.Cells.Item(29, 1).Value = 'knee pads weightlifting'
$ws.Cells.Item(30, 1).Value = 'knee pads youth'
$ws.Cells.Item(31, 1).Value = 'knee pants for men'
$ws.Cells.Item(32, 1).Value = 'knee protection'
$ws.Cells.Item(33, 1).Value = 'knee protector gym'
$ws.Cells.Item(34, 1).Value = 'knee protectors'
$ws.Cells.Item(35, 1).Value = 'knee replacement aids'
$ws.Cells.Item(36, 1).Value = 'knee support gym'
$ws.Cells.Item(37, 1).Value = 'knee support pants'
$ws.Cells.Item(38, 1).Value = 'knee support pants men'
$ws.Cells.Item(39, 1).Value = 'kneepads for volleyball'
$ws.Cells.Item(40, 1).Value = 'leg compression pants'
$ws.Cells.Item(41, 1).Value = 'leg compression tights'
$ws.Cells.Item(42, 1).Value = 'leg protectors for men'
$ws.Cells.Item(43, 1).Value = 'legging for basketball boys'
$ws.Cells.Item(44, 1).Value = 'leggings for men sport gym'
$ws.Cells.Item(45, 1).Value = 'leggings knee'
$ws.Cells.Item(46, 1).Value = 'leggings medium'
$ws.Cells.Item(47, 1).Value = 'leggings mens'
$ws.Cells.Item(48, 1).Value = 'leggins training'
$ws.Cells.Item(49, 1).Value = 'lightweight athletic pants for men'
$ws.Cells.Item(50, 1).Value = 'lightweight pants men'
$ws.Cells.Item(51, 1).Value = 'lightweight sports pants men'
$ws.Cells.Item(52, 1).Value = 'mcdavid basketball knee pads 6446'
$ws.Cells.Item(53, 1).Value = 'mcdavid basketball knee pads black'
$ws.Cells.Item(54, 1).Value = 'men basketball tights'
$ws.Cells.Item(55, 1).Value = 'men capri pants'
$ws.Cells.Item(56, 1).Value = 'men compression pants'
$ws.Cells.Item(57, 1).Value = 'men leggings pack'
$ws.Cells.Item(58, 1).Value = 'men leggings running'
$ws.Cells.Item(59, 1).Value = 'men leggings tall'
$ws.Cells.Item(60, 1).Value = 'men running tights nike'
$ws.Cells.Item(61, 1).Value = 'men tights legging'
$ws.Cells.Item(62, 1).Value = 'mens athletic compression pants'
$ws.Cells.Item(63, 1).Value = 'mens athletic leggings black'
$ws.Cells.Item(64, 1).Value = 'mens basketball knee pads'
$ws.Cells.Item(65, 1).Value = 'mens compression 3 4 pants'
$ws.Cells.Item(66, 1).Value = 'mens compression leggings 3 4'
$ws.Cells.Item(67, 1).Value = 'mens compression tights pants'
$ws.Cells.Item(68, 1).Value = 'mens cycling tights'
$ws.Cells.Item(69, 1).Value = 'mens hiking pants lightweight'
$ws.Cells.Item(70, 1).Value = 'mens jogging tights'
$ws.Cells.Item(71, 1).Value = 'mens lacrosse pads'
$ws.Cells.Item(72, 1).Value = 'mens leggings'
$ws.Cells.Item(73, 1).Value = 'mens lightweight workout pants'
$ws.Cells.Item(74, 1).Value = 'mens running tights green'
$ws.Cells.Item(75, 1).Value = 'mens running tights orange'
$ws.Cells.Item(76, 1).Value = 'mens running tights yellow'
$ws.Cells.Item(77, 1).Value = 'mens soccer clothing'
$ws.Cells.Item(78, 1).Value = 'mens stretch thermal pants'
$ws.Cells.Item(79, 1).Value = 'mens swim leggings'
$ws.Cells.Item(80, 1).Value = 'mens tights basketball'
$ws.Cells.Item(81, 1).Value = 'mens tights capri'
$ws.Cells.Item(82, 1).Value = 'mens tights leggings'
$ws.Cells.Item(83, 1).Value = 'mens underarmour snow pants'
$ws.Cells.Item(84, 1).Value = 'mens workout pants'
$ws.Cells.Item(85, 1).Value = 'multicam pants with knee pads'
$ws.Cells.Item(86, 1).Value = 'nike basketball pads'
$ws.Cells.Item(87, 1).Value = 'nike volleyball knee pads youth girls'
$ws.Cells.Item(88, 1).Value = 'padded compression'
$ws.Cells.Item(89, 1).Value = 'padded soccer pants'
$ws.Cells.Item(90, 1).Value = 'pain in thigh joint'
$ws.Cells.Item(91, 1).Value = 'pant knee pad inserts'
$ws.Cells.Item(92, 1).Value = 'pantalon con rodilleras'
$ws.Cells.Item(93, 1).Value = 'pants with knee pads'
$ws.Cells.Item(94, 1).Value = 'polyester pants men pants'
$ws.Cells.Item(95, 1).Value = 'protective knee pads for men'
$ws.Cells.Item(96, 1).Value = 'rash guard men bjj'
$ws.Cells.Item(97, 1).Value = 'reebok compression pants'
$ws.Cells.Item(98, 1).Value = 'running compression leg'
$ws.Cells.Item(99, 1).Value = 'running knee'
$ws.Cells.Item(100, 1).Value = 'running pants youth'
